# Adds a new "localdb" command-type column to the '#system' sheet and wires
# up a matching named range, per the commit:
#   "NEW command type to create, maintain and manipulate a local-only
#    relational database ..."
#
# Net effect on the '#system' sheet:
#   - a new column is inserted before column N; everything from N..AC
#     (mail, number, pdf, rdbms, redis, sms, sound, ssh, step, web,
#     webalert, webcookie, ws, ws.async, xml, macro, ...) shifts one
#     column to the right (O..AD).
#   - the freed-up column N is populated with the new "localdb" function
#     list (header + 6 functions).
#   - the "target" category list in column A gets a new "localdb" entry
#     inserted alphabetically (between "json" and "macro"), so rows
#     14..29 shift down to 15..30.
#   - defined names are updated to reflect the new column letters, and a
#     new "localdb" defined name is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a new column before N - shifts N:AC to O:AD (values, not
#    just formatting), freeing up column N for the new "localdb" list.
# ---------------------------------------------------------------------
$ws.Columns("N").Insert()

# ---------------------------------------------------------------------
# 2) Insert the new "localdb" entry into the "target" category list
#    (column A), shifting rows 14:29 down to 15:30. Range.Insert() in
#    this host operates at row granularity (it would drag every column
#    along), so shift column A manually, cell by cell, bottom-up.
# ---------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(14, 1).Value = "localdb"

# ---------------------------------------------------------------------
# 3) Populate the new column N with the "localdb" header + functions.
# ---------------------------------------------------------------------
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 4) Fix up the defined names that pointed at columns N..AC (now shifted
#    one column right), the "target" list (now one row taller), and add
#    the new "localdb" defined name.
# ---------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo = "='#system'!`$O`$2:`$O`$4"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
